$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3909.0833
$ws.Range("I40").Value = 2858.1667
$ws.Range("J40").Value = 4960
$ws.Range("K40").Value = 2858.1667
$ws.Range("L40").Value = 4960
$ws.Range("M40").Value = -2683.1667
$ws.Range("N40").Value = -5310
$ws.Range("H92").Value = 55555924
$ws.Range("I92").Value = 83333490
$ws.Range("J92").Value = 798.5
$ws.Range("K92").Value = 83333490
$ws.Range("L92").Value = 798.5
$ws.Range("M92").Value = -83332242
$ws.Range("H97").Value = 1880.125
$ws.Range("J97").Value = 1945.4667
$ws.Range("L97").Value = 5836.4001
$ws.Range("N97").Value = -6828.4001
$ws.Range("H100").Value = 1645.0588
$ws.Range("I100").Value = 1241.0625
$ws.Range("J100").Value = 2004.1666
$ws.Range("K100").Value = 1241.0625
$ws.Range("L100").Value = 2004.1666
$ws.Range("M100").Value = -700.0625
$ws.Range("N100").Value = -3086.1666
$ws.Range("H112").Value = 2156.9736
$ws.Range("J112").Value = 1859.7273
$ws.Range("L112").Value = 5579.1819
$ws.Range("N112").Value = -7795.1819
$ws.Range("H137").Value = 2888.2134
$ws.Range("I137").Value = 3029.647
$ws.Range("K137").Value = 9088.940999999999
$ws.Range("M137").Value = -6538.940999999999
$ws.Range("H139").Value = 69948.664
$ws.Range("J139").Value = 69948.664
$ws.Range("L139").Value = 69948.664
$ws.Range("N139").Value = -80228.664
$ws.Range("H141").Value = 3478.923
$ws.Range("I141").Value = 3478.923
$ws.Range("K141").Value = 10436.769
$ws.Range("M141").Value = -5256.769
$ws.Range("N92").Value = -3294.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3914.2185
$ws.Range("I32").Value = 3970.7073
$ws.Range("J32").Value = 2987.8
$ws.Range("K32").Value = 3970.7073
$ws.Range("L32").Value = 2987.8
$ws.Range("M32").Value = -3683.7073
$ws.Range("N32").Value = -3561.8
$ws.Range("H63").Value = 2139.4
$ws.Range("I63").Value = 2139.4
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2139.4
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1453.4
$ws.Range("H66").Value = 2139.4
$ws.Range("I66").Value = 2139.4
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 10697
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -7265
$ws.Range("H74").Value = 1801
$ws.Range("I74").Value = 1801
$ws.Range("K74").Value = 1801
$ws.Range("M74").Value = -927
$ws.Range("H77").Value = 1801
$ws.Range("I77").Value = 1801
$ws.Range("K77").Value = 9005
$ws.Range("M77").Value = -4637
$ws.Range("H125").Value = 60379.832
$ws.Range("J125").Value = 60379.832
$ws.Range("L125").Value = 60379.832
$ws.Range("N125").Value = -70219.83199999999
$ws.Range("H132").Value = 1491.1052
$ws.Range("I132").Value = 1491.1052
$ws.Range("K132").Value = 4473.3156
$ws.Range("M132").Value = -1943.3156
$ws.Range("H133").Value = 18285.715
$ws.Range("J133").Value = 18285.715
$ws.Range("L133").Value = 18285.715
$ws.Range("N133").Value = -23345.715
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 118499
$ws.Range("J55").Value = 118499
$ws.Range("L55").Value = 118499
$ws.Range("N55").Value = -119045
$ws.Range("H99").Value = 1393.0741
$ws.Range("I99").Value = 1353.9474
$ws.Range("J99").Value = 1486
$ws.Range("K99").Value = 1353.9474
$ws.Range("L99").Value = 1486
$ws.Range("M99").Value = 144.0526
$ws.Range("N99").Value = -4482

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4736.8066
$ws.Range("I16").Value = 3880.7222
$ws.Range("J16").Value = 5922.154
$ws.Range("K16").Value = 3880.7222
$ws.Range("L16").Value = 5922.154
$ws.Range("M16").Value = -3593.7222
$ws.Range("N16").Value = -6496.154
$ws.Range("H31").Value = 2633.0833
$ws.Range("I31").Value = 2425.3447
$ws.Range("J31").Value = 3493.7144
$ws.Range("K31").Value = 2425.3447
$ws.Range("L31").Value = 3493.7144
$ws.Range("M31").Value = -2130.3447
$ws.Range("N31").Value = -4083.7144
$ws.Range("H34").Value = 2633.0833
$ws.Range("I34").Value = 2425.3447
$ws.Range("J34").Value = 3493.7144
$ws.Range("K34").Value = 2425.3447
$ws.Range("L34").Value = 3493.7144
$ws.Range("M34").Value = -2223.3447
$ws.Range("N34").Value = -3897.7144
$ws.Range("H64").Value = 57307.617
$ws.Range("J64").Value = 57307.617
$ws.Range("L64").Value = 57307.617
$ws.Range("N64").Value = -57803.617
$ws.Range("H67").Value = 57307.617
$ws.Range("J67").Value = 57307.617
$ws.Range("L67").Value = 57307.617
$ws.Range("N67").Value = -59023.617
$ws.Range("H113").Value = 4736.8066
$ws.Range("I113").Value = 3880.7222
$ws.Range("J113").Value = 5922.154
$ws.Range("K113").Value = 3880.7222
$ws.Range("L113").Value = 5922.154
$ws.Range("M113").Value = -1710.7222
$ws.Range("N113").Value = -10262.154
$ws.Range("H124").Value = 59040
$ws.Range("J124").Value = 59040
$ws.Range("L124").Value = 59040
$ws.Range("N124").Value = -63950
$ws.Range("H134").Value = 1533.6567
$ws.Range("I134").Value = 1571.0566
$ws.Range("J134").Value = 1392.0714
$ws.Range("K134").Value = 4713.1698
$ws.Range("L134").Value = 4176.2142
$ws.Range("M134").Value = -2178.1698
$ws.Range("N134").Value = -9246.2142

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 15712.611
$ws.Range("I56").Value = 15712.611
$ws.Range("K56").Value = 15712.611
$ws.Range("M56").Value = -15182.611
$ws.Range("H97").Value = 806.1539
$ws.Range("I97").Value = 709.5
$ws.Range("K97").Value = 2128.5
$ws.Range("M97").Value = -1632.5
$ws.Range("H113").Value = 810.53845
$ws.Range("I113").Value = 550.1875
$ws.Range("J113").Value = 1227.1
$ws.Range("K113").Value = 1650.5625
$ws.Range("L113").Value = 3681.3
$ws.Range("M113").Value = 519.4375
$ws.Range("N113").Value = -8021.299999999999
$ws.Range("H122").Value = 1730.2941
$ws.Range("I122").Value = 589
$ws.Range("J122").Value = 1801.625
$ws.Range("K122").Value = 5301
$ws.Range("L122").Value = 16214.625
$ws.Range("M122").Value = -2851
$ws.Range("N122").Value = -21114.625
$ws.Range("H132").Value = 2401.2964
$ws.Range("J132").Value = 3219.6667
$ws.Range("L132").Value = 28977.0003
$ws.Range("N132").Value = -34037.0003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20572.428
$ws.Range("I70").Value = 52004
$ws.Range("J70").Value = 7999.8
$ws.Range("K70").Value = 52004
$ws.Range("L70").Value = 7999.8
$ws.Range("M70").Value = -51734
$ws.Range("N70").Value = -8539.799999999999
$ws.Range("H73").Value = 20572.428
$ws.Range("I73").Value = 52004
$ws.Range("J73").Value = 7999.8
$ws.Range("K73").Value = 52004
$ws.Range("L73").Value = 7999.8
$ws.Range("M73").Value = -51068
$ws.Range("N73").Value = -9871.799999999999
$ws.Range("H102").Value = 4241.7173
$ws.Range("I102").Value = 3358.2122
$ws.Range("J102").Value = 6484.4614
$ws.Range("K102").Value = 3358.2122
$ws.Range("L102").Value = 6484.4614
$ws.Range("M102").Value = -1736.2122
$ws.Range("N102").Value = -9728.4614
$ws.Range("H126").Value = 6641.9
$ws.Range("I126").Value = 7374.5
$ws.Range("J126").Value = 6153.5
$ws.Range("K126").Value = 22123.5
$ws.Range("L126").Value = 18460.5
$ws.Range("M126").Value = -19653.5
$ws.Range("N126").Value = -23400.5
$ws.Range("H132").Value = 5193.5117
$ws.Range("I132").Value = 5020.073
$ws.Range("K132").Value = 15060.219
$ws.Range("M132").Value = -12530.219

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 25924.6
$ws.Range("I132").Value = 29277.348
$ws.Range("J132").Value = 14908.429
$ws.Range("K132").Value = 87832.04400000001
$ws.Range("L132").Value = 44725.287
$ws.Range("M132").Value = -85302.04400000001
$ws.Range("N132").Value = -49785.287

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4562.533
$ws.Range("I96").Value = 4371.778
$ws.Range("J96").Value = 4848.6665
$ws.Range("K96").Value = 4371.778
$ws.Range("L96").Value = 4848.6665
$ws.Range("M96").Value = -2998.778
$ws.Range("N96").Value = -7594.6665
$ws.Range("H132").Value = 2129.3333
$ws.Range("I132").Value = 1497.2745
$ws.Range("J132").Value = 3920.1667
$ws.Range("K132").Value = 4491.8235
$ws.Range("L132").Value = 11760.5001
$ws.Range("M132").Value = -1961.8235
$ws.Range("N132").Value = -16820.5001
$ws.Range("H136").Value = 2109.9358
$ws.Range("I136").Value = 2097.096
$ws.Range("K136").Value = 6291.288
$ws.Range("M136").Value = -3741.288
